$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H (HALKBANK) values that were previously blank inline-string cells
$ws.Range("H3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("H8").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("H9").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("H10").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("H11").Value = "3,05 TL - 6,1 TL - 76,18 TL"
$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"
$ws.Range("H14").Value = "2.100 TL - 4.300 TL"

# Updated figures for AKBANK (D13) and ISBANKASI (E13)
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
